# Remove every <w:contextualSpacing w:val="0"/> paragraph property from
# every paragraph in the document. This property was dropped from every
# paragraph's <w:pPr> in the target revision; Word's exposed object model
# in this runtime has no dedicated ContextualSpacing property, so we
# surgically edit the paragraph's underlying WordprocessingML via
# Range.XML()/Range.InsertXML (both standard Word COM members).

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $range = $para.Range

    $full = $range.XML()

    if ($full -notmatch '(?s)<w:body>(.*)</w:body>') {
        continue
    }
    $bodyInner = $Matches[1]

    if ($bodyInner -notmatch '(?s)^(<w:p\b.*?</w:p>)') {
        continue
    }
    $paraXml = $Matches[1]

    if ($paraXml -notmatch '<w:contextualSpacing\b[^/]*/>') {
        continue
    }

    $newXml = $paraXml -replace '<w:contextualSpacing\b[^/]*/>', ''
    $range.InsertXML($newXml)
}

Write-Output "done: processed $count paragraphs"
